# Add new memes/stickers to the meme_bot_db workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (Meme, StickerID) pairs to append right after the existing data
# (last populated row is 158).
$memes = @(
    "la falta de respeto abunda por aquí | falta de respeto",
    "fuera impulso de idiotez | impulso de idiotez | no fuera impulso de idiotez",
    "Cuantas veces tenemos que darle una leccion anciano | darle una leccion",
    "alan que pendejada hiciste | que hiciste",
    "tuve ese sueño de nuevo | gary tuve ese sueño de nuevo",
    "mi lente de contacto | oh mi lente de contacto",
    "deja de hablar y pegale | pegale",
    "te parece que somos ricos",
    "demasiada comedia | no soportamos tanta comedia"
)

$stickers = @(
    "CAACAgEAAxkBAAI7EmB_Jg7ylTo7OuX768XDjWRZq8DaAAKEAQACFhT4R8OrJw7LLS6ZHwQ",
    "CAACAgEAAxkBAAI7FGB_JjH5EvUZGkrZBNG1NlDMjNE8AAISAQACHi75R5k_Ay8pdjtPHwQ",
    "CAACAgEAAxkBAAI7FmB_Jl1X_ZyAfmTywPCy1Zs1jIhRAALUAQACsUD5R0xHerSaTa56HwQ",
    "CAACAgEAAxkBAAI7GGB_JoACCKnopQuVVd4xONC7ipWlAAKPAQAC4c34R05QjFIOVXvbHwQ",
    "CAACAgEAAxkBAAI7GmB_JqAMdkcyGpILwieGLlimugIgAAKXAQACpAH5R8QXWC-dXoZgHwQ",
    "CAACAgEAAxkBAAI7HGB_JtlvA7p8p6A3B2CBKguDsPjUAALxAQACXEP5R7ZWljylnKu-HwQ",
    "CAACAgEAAxkBAAI7HmB_JvEX1Fsld4tb8vDccdGJbbBhAAJXAQACWM0BRLaLkPndZ_QOHwQ",
    "CAACAgEAAxkBAAI7IGB_JwJ6AAFj7OQlovxDSAGc0iICEQACMgEAAqKb-EfauivMw5lHlR8E",
    "CAACAgEAAxkBAAI7ImB_Jw_AjGb1Btul7XP2rgABCgVxWAACVQEAAvm6-EdKo9efOIHhJx8E"
)

$startRow = 159
$count = $memes.Length

# Copy the formatting of the last populated row (158) down onto the new
# rows first, so the new cells pick up the same style (s="4") as row 158.
$endRow = $startRow + $count - 1
$ws.Range("A158:B158").Copy()
$ws.Range("A" + $startRow + ":B" + $endRow).PasteSpecial(-4122)

for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $memes[$i]
    $ws.Cells.Item($r, 2).Value = $stickers[$i]
}

# The trailing blank row 993 was removed from the sheet.
$ws.Rows.Item(993).Delete()
